$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1091.8158
$ws.Cells.Item(70, 9).Value = 999.8095
$ws.Cells.Item(70, 10).Value = 1205.4706
$ws.Cells.Item(70, 11).Value = 2999.4285
$ws.Cells.Item(70, 12).Value = 3616.4118
$ws.Cells.Item(70, 13).Value = -2729.4285
$ws.Cells.Item(70, 14).Value = -4156.4118
$ws.Cells.Item(73, 8).Value = 1091.8158
$ws.Cells.Item(73, 9).Value = 999.8095
$ws.Cells.Item(73, 10).Value = 1205.4706
$ws.Cells.Item(73, 11).Value = 2999.4285
$ws.Cells.Item(73, 12).Value = 3616.4118
$ws.Cells.Item(73, 13).Value = -2063.4285
$ws.Cells.Item(73, 14).Value = -5488.4118
$ws.Cells.Item(96, 8).Value = 1177.8695
$ws.Cells.Item(96, 9).Value = 813.4666999999999
$ws.Cells.Item(96, 10).Value = 1861.125
$ws.Cells.Item(96, 11).Value = 2440.4001
$ws.Cells.Item(96, 12).Value = 5583.375
$ws.Cells.Item(96, 13).Value = -1067.4001
$ws.Cells.Item(96, 14).Value = -8329.375
$ws.Cells.Item(100, 8).Value = 2564.2307
$ws.Cells.Item(100, 9).Value = 1938.75
$ws.Cells.Item(100, 10).Value = 2842.2222
$ws.Cells.Item(100, 11).Value = 1938.75
$ws.Cells.Item(100, 12).Value = 2842.2222
$ws.Cells.Item(100, 13).Value = -1397.75
$ws.Cells.Item(100, 14).Value = -3924.2222
$ws.Cells.Item(125, 8).Value = 252725
$ws.Cells.Item(125, 9).Value = 1000000
$ws.Cells.Item(125, 10).Value = 3633.3333
$ws.Cells.Item(125, 11).Value = 9000000
$ws.Cells.Item(125, 12).Value = 32699.9997
$ws.Cells.Item(125, 13).Value = -8997540
$ws.Cells.Item(125, 14).Value = -37619.9997
$ws.Cells.Item(127, 8).Value = 597.2
$ws.Cells.Item(127, 9).Value = 519.8461
$ws.Cells.Item(127, 11).Value = 1559.5383
$ws.Cells.Item(127, 13).Value = 3400.4617
$ws.Cells.Item(132, 8).Value = 5255.431
$ws.Cells.Item(132, 9).Value = 4215.089
$ws.Cells.Item(132, 10).Value = 8856.615
$ws.Cells.Item(132, 11).Value = 12645.267
$ws.Cells.Item(132, 12).Value = 26569.845
$ws.Cells.Item(132, 13).Value = -10115.267
$ws.Cells.Item(132, 14).Value = -31629.845
$ws.Cells.Item(137, 8).Value = 2431.6538
$ws.Cells.Item(137, 9).Value = 3038.2666
$ws.Cells.Item(137, 10).Value = 1604.4546
$ws.Cells.Item(137, 11).Value = 9114.799800000001
$ws.Cells.Item(137, 12).Value = 4813.3638
$ws.Cells.Item(137, 13).Value = -6564.799800000001
$ws.Cells.Item(137, 14).Value = -9913.363799999999
$ws.Cells.Item(138, 8).Value = 1831.674
$ws.Cells.Item(138, 9).Value = 1323.6451
$ws.Cells.Item(138, 10).Value = 2881.6
$ws.Cells.Item(138, 11).Value = 3970.9353
$ws.Cells.Item(138, 12).Value = 8644.799999999999
$ws.Cells.Item(138, 13).Value = 1169.0647
$ws.Cells.Item(138, 14).Value = -18924.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6468.891
$ws.Cells.Item(32, 9).Value = 6373.1133
$ws.Cells.Item(32, 10).Value = 9007
$ws.Cells.Item(32, 11).Value = 6373.1133
$ws.Cells.Item(32, 12).Value = 9007
$ws.Cells.Item(32, 13).Value = -6086.1133
$ws.Cells.Item(32, 14).Value = -9581
$ws.Cells.Item(74, 8).Value = 1945.52
$ws.Cells.Item(74, 9).Value = 1995.8667
$ws.Cells.Item(74, 10).Value = 1870
$ws.Cells.Item(74, 11).Value = 1995.8667
$ws.Cells.Item(74, 12).Value = 1870
$ws.Cells.Item(74, 13).Value = -1121.8667
$ws.Cells.Item(74, 14).Value = -3618
$ws.Cells.Item(77, 8).Value = 1945.52
$ws.Cells.Item(77, 9).Value = 1995.8667
$ws.Cells.Item(77, 10).Value = 1870
$ws.Cells.Item(77, 11).Value = 9979.333500000001
$ws.Cells.Item(77, 12).Value = 9350
$ws.Cells.Item(77, 13).Value = -5611.333500000001
$ws.Cells.Item(77, 14).Value = -18086
$ws.Cells.Item(96, 8).Value = 26192
$ws.Cells.Item(96, 10).Value = 26192
$ws.Cells.Item(96, 12).Value = 26192
$ws.Cells.Item(96, 14).Value = -31684
$ws.Cells.Item(102, 8).Value = 1474.1666
$ws.Cells.Item(102, 9).Value = 1523.3334
$ws.Cells.Item(102, 10).Value = 1326.6666
$ws.Cells.Item(102, 11).Value = 1523.3334
$ws.Cells.Item(102, 12).Value = 1326.6666
$ws.Cells.Item(102, 13).Value = 98.66660000000002
$ws.Cells.Item(102, 14).Value = -4570.6666
$ws.Cells.Item(110, 8).Value = 2223.8333
$ws.Cells.Item(110, 9).Value = 1067.2667
$ws.Cells.Item(110, 10).Value = 4151.4443
$ws.Cells.Item(110, 11).Value = 1067.2667
$ws.Cells.Item(110, 12).Value = 4151.4443
$ws.Cells.Item(110, 13).Value = 977.7333000000001
$ws.Cells.Item(110, 14).Value = -8241.444299999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1004.3333
$ws.Cells.Item(107, 9).Value = 1000
$ws.Cells.Item(107, 10).Value = 1006.5
$ws.Cells.Item(107, 11).Value = 1000
$ws.Cells.Item(107, 12).Value = 1006.5
$ws.Cells.Item(107, 13).Value = 920
$ws.Cells.Item(107, 14).Value = -4846.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2936
$ws.Cells.Item(31, 9).Value = 2325.8948
$ws.Cells.Item(31, 10).Value = 6800
$ws.Cells.Item(31, 11).Value = 2325.8948
$ws.Cells.Item(31, 12).Value = 6800
$ws.Cells.Item(31, 13).Value = -2030.8948
$ws.Cells.Item(31, 14).Value = -7390
$ws.Cells.Item(34, 8).Value = 2936
$ws.Cells.Item(34, 9).Value = 2325.8948
$ws.Cells.Item(34, 10).Value = 6800
$ws.Cells.Item(34, 11).Value = 2325.8948
$ws.Cells.Item(34, 12).Value = 6800
$ws.Cells.Item(34, 13).Value = -2123.8948
$ws.Cells.Item(34, 14).Value = -7204
$ws.Cells.Item(58, 8).Value = 4388.5347
$ws.Cells.Item(58, 9).Value = 1606.36
$ws.Cells.Item(58, 10).Value = 8252.666999999999
$ws.Cells.Item(58, 11).Value = 1606.36
$ws.Cells.Item(58, 12).Value = 8252.666999999999
$ws.Cells.Item(58, 13).Value = -1403.36
$ws.Cells.Item(58, 14).Value = -8658.666999999999
$ws.Cells.Item(136, 8).Value = 4388.5347
$ws.Cells.Item(136, 9).Value = 1606.36
$ws.Cells.Item(136, 10).Value = 8252.666999999999
$ws.Cells.Item(136, 11).Value = 4819.08
$ws.Cells.Item(136, 12).Value = 24758.001
$ws.Cells.Item(136, 13).Value = -2269.08
$ws.Cells.Item(136, 14).Value = -29858.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 743.9091
$ws.Cells.Item(113, 9).Value = 602
$ws.Cells.Item(113, 10).Value = 797.125
$ws.Cells.Item(113, 11).Value = 1806
$ws.Cells.Item(113, 12).Value = 2391.375
$ws.Cells.Item(113, 13).Value = 364
$ws.Cells.Item(113, 14).Value = -6731.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2233.3333
$ws.Cells.Item(113, 9).Value = 2111.111
$ws.Cells.Item(113, 10).Value = 2600
$ws.Cells.Item(113, 11).Value = 2111.111
$ws.Cells.Item(113, 12).Value = 2600
$ws.Cells.Item(113, 13).Value = 58.88900000000012
$ws.Cells.Item(113, 14).Value = -6940
$ws.Cells.Item(122, 8).Value = 717899.0600000001
$ws.Cells.Item(122, 9).Value = 1002860.7
$ws.Cells.Item(122, 10).Value = 5495
$ws.Cells.Item(122, 11).Value = 3008582.1
$ws.Cells.Item(122, 12).Value = 16485
$ws.Cells.Item(122, 13).Value = -3006132.1
$ws.Cells.Item(122, 14).Value = -21385
$ws.Cells.Item(126, 8).Value = 3417.2173
$ws.Cells.Item(126, 9).Value = 3739.2307
$ws.Cells.Item(126, 10).Value = 2998.6
$ws.Cells.Item(126, 11).Value = 11217.6921
$ws.Cells.Item(126, 12).Value = 8995.799999999999
$ws.Cells.Item(126, 13).Value = -8747.6921
$ws.Cells.Item(126, 14).Value = -13935.8
$ws.Cells.Item(132, 8).Value = 2842.75
$ws.Cells.Item(132, 9).Value = 4064
$ws.Cells.Item(132, 10).Value = 2560.923
$ws.Cells.Item(132, 11).Value = 12192
$ws.Cells.Item(132, 12).Value = 7682.768999999999
$ws.Cells.Item(132, 13).Value = -9662
$ws.Cells.Item(132, 14).Value = -12742.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1567.2693
$ws.Cells.Item(7, 9).Value = 1669.7333
$ws.Cells.Item(7, 10).Value = 1427.5454
$ws.Cells.Item(7, 11).Value = 1669.7333
$ws.Cells.Item(7, 12).Value = 1427.5454
$ws.Cells.Item(7, 13).Value = -1557.7333
$ws.Cells.Item(7, 14).Value = -1651.5454
$ws.Cells.Item(16, 8).Value = 1344.8
$ws.Cells.Item(16, 9).Value = 1312.7693
$ws.Cells.Item(16, 10).Value = 1553
$ws.Cells.Item(16, 11).Value = 1312.7693
$ws.Cells.Item(16, 12).Value = 1553
$ws.Cells.Item(16, 13).Value = -1142.7693
$ws.Cells.Item(16, 14).Value = -1893
$ws.Cells.Item(40, 8).Value = 45457240
$ws.Cells.Item(40, 9).Value = 76924904
$ws.Cells.Item(40, 10).Value = 3945.5557
$ws.Cells.Item(40, 11).Value = 76924904
$ws.Cells.Item(40, 12).Value = 3945.5557
$ws.Cells.Item(40, 13).Value = -76924768
$ws.Cells.Item(40, 14).Value = -4217.5557
$ws.Cells.Item(95, 8).Value = 13490.667
$ws.Cells.Item(95, 10).Value = 13490.667
$ws.Cells.Item(95, 12).Value = 13490.667
$ws.Cells.Item(95, 14).Value = -18982.667
$ws.Cells.Item(100, 8).Value = 2798.9092
$ws.Cells.Item(100, 9).Value = 2220
$ws.Cells.Item(100, 10).Value = 3129.7144
$ws.Cells.Item(100, 11).Value = 2220
$ws.Cells.Item(100, 12).Value = 3129.7144
$ws.Cells.Item(100, 13).Value = -1679
$ws.Cells.Item(100, 14).Value = -4211.7144
$ws.Cells.Item(126, 8).Value = 1567.2693
$ws.Cells.Item(126, 9).Value = 1669.7333
$ws.Cells.Item(126, 10).Value = 1427.5454
$ws.Cells.Item(126, 11).Value = 5009.199900000001
$ws.Cells.Item(126, 12).Value = 4282.6362
$ws.Cells.Item(126, 13).Value = -2539.199900000001
$ws.Cells.Item(126, 14).Value = -9222.636200000001
$ws.Cells.Item(132, 8).Value = 29414638
$ws.Cells.Item(132, 9).Value = 45457210
$ws.Cells.Item(132, 10).Value = 3257.6667
$ws.Cells.Item(132, 11).Value = 136371630
$ws.Cells.Item(132, 12).Value = 9773.000100000001
$ws.Cells.Item(132, 13).Value = -136369100
$ws.Cells.Item(132, 14).Value = -14833.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).Value = $null
$ws.Cells.Item(126, 8).Value = 8102.421
$ws.Cells.Item(126, 9).Value = 8385.888999999999
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 25157.667
$ws.Cells.Item(126, 12).Value = 9000
$ws.Cells.Item(126, 13).Value = -22687.667
$ws.Cells.Item(126, 14).Value = -13940
$ws.Cells.Item(132, 8).Value = 3265.8096
$ws.Cells.Item(132, 9).Value = 3664
$ws.Cells.Item(132, 10).Value = 3106.5334
$ws.Cells.Item(132, 11).Value = 10992
$ws.Cells.Item(132, 12).Value = 9319.600199999999
$ws.Cells.Item(132, 13).Value = -8462
$ws.Cells.Item(132, 14).Value = -14379.6002
